$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 444
$ws.Range("B5").Value = "NewBorn"

$ws.Range("A9").Select()
